# Save before changing branches
# Update simulation output values on Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 558
$ws.Range("G2").Value = 186

$ws.Range("F4").Value = 2444
$ws.Range("G4").Value = 1540
$ws.Range("I4").Value = 34.645338449995506

$ws.Range("F5").Value = 4979
$ws.Range("G5").Value = 3999
$ws.Range("I5").Value = 45.152801868133103

$ws.Range("F6").Value = 6409
$ws.Range("G6").Value = 4330
$ws.Range("I6").Value = 47.598651780239997

$ws.Range("F7").Value = 420
$ws.Range("G7").Value = 21

$ws.Range("F8").Value = 1296
$ws.Range("G8").Value = 408

$ws.Range("F9").Value = 3481
$ws.Range("G9").Value = 2369
$ws.Range("I9").Value = 45.06247185710609

$ws.Range("F10").Value = 4387
$ws.Range("G10").Value = 2556
$ws.Range("I10").Value = 49.849273946872302

$ws.Range("F11").Value = 6312
$ws.Range("G11").Value = 4604
$ws.Range("I11").Value = 49.187662705293334

$ws.Range("F12").Value = 431
$ws.Range("G12").Value = 20

$ws.Range("F13").Value = 990
$ws.Range("G13").Value = 205

$ws.Range("F14").Value = 2639
$ws.Range("G14").Value = 1369
$ws.Range("I14").Value = 42.416029932403752

$ws.Range("F15").Value = 4648
$ws.Range("G15").Value = 3108
$ws.Range("I15").Value = 48.989936596123023

$ws.Range("B16").Value = 0
$ws.Range("F16").Value = 7860
$ws.Range("G16").Value = 6400
$ws.Range("I16").Value = 50.458955984150471
